$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39/40 swap: ARBITRUM <-> Stellar (coin, link, price, volume) ---
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.114"
$ws.Range("E39").Value = "  -1.32%  "

$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.81"
$ws.Range("E40").Value = "  -2.43%  "

# --- Price (D) / Volume(1h) (E) updates for remaining rows ---
$ws.Range("D2").Value = "42.077.37"
$ws.Range("E2").Value = "  -0.80%  "
$ws.Range("D3").Value = "2.259.19"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.13"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.68"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  -1.01%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.487"
$ws.Range("E9").Value = "  -1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.89"
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0785"
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.81"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("D14").Value = "2.613.40"
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.63"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "2.255.48"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.786"
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").Value = "41.939.42"
$ws.Range("E18").Value = "  -0.96%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.21"
$ws.Range("E19").Value = "  -3.84%  "
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.97"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.46"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.10"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  -1.68%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.41"
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.82"
$ws.Range("E28").Value = "  +2.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.52"
$ws.Range("E29").Value = "  -0.72%  "
$ws.Range("E30").Value = "  +0.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.24"
$ws.Range("E31").Value = "  +1.87%  "
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.13"
$ws.Range("E34").Value = "  +1.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.60"
$ws.Range("E35").Value = "  +2.24%  "
$ws.Range("E36").Value = "  -3.05%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("E38").Value = "  -5.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.08"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.32"
$ws.Range("E42").Value = "  +1.48%  "
$ws.Range("D43").Value = "1.947.11"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.97"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("E45").Value = "  -2.19%  "
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.88"
$ws.Range("E47").Value = "  -4.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.18"
$ws.Range("E48").Value = "  -1.01%  "
$ws.Range("D49").Value = "2.485.47"
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.77"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.97"
$ws.Range("E51").Value = "  -0.21%  "
